$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ITI")

$ws.Range("J2").Value = "Resolvido"
$ws.Range("J3").Value = "Resolvido"
$ws.Range("J4").Value = "Resolvido"
$ws.Range("J5").Value = "Resolvido"
